# edit.ps1 - applies the "Updated cryptos list" GitHub Actions refresh
# (price + 1h-volume updates, plus a VeChain/Cronos rank swap at rows 43-44)
# to the cryptos worksheet, cell by cell, matching the target OOXML diff.
#
# Column D ("Price") holds values that are stored as *text* in the workbook
# (e.g. "37.208.33", "54.90") even though many of them look numeric. Plain
# COM assignment of a numeric-looking string (e.g. "54.45") would be
# auto-converted to a real number by Excel, losing the text type and any
# trailing zeros. To keep those cells as text we prefix the literal with a
# leading apostrophe ('), exactly as a user would do when typing into Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '37.035.14'
$ws.Range('E2').Value = '  -1.04%  '

# Row 3
$ws.Range('D3').Value = '1.979.17'
$ws.Range('E3').Value = '  -2.95%  '

# Row 4
$ws.Range('E4').Value = '  +0.26%  '

# Row 5
$ws.Range('D5').Value = '''250.08'
$ws.Range('E5').Value = '  -3.23%  '

# Row 6
$ws.Range('D6').Value = '''0.602'
$ws.Range('E6').Value = '  -3.45%  '

# Row 7
$ws.Range('E7').Value = '  +0.35%  '

# Row 8
$ws.Range('D8').Value = '''54.45'
$ws.Range('E8').Value = '  -6.73%  '

# Row 9
$ws.Range('D9').Value = '''0.372'
$ws.Range('E9').Value = '  -4.80%  '

# Row 10
$ws.Range('D10').Value = '''0.0750'
$ws.Range('E10').Value = '  -7.50%  '

# Row 11
$ws.Range('D11').Value = '''0.0994'
$ws.Range('E11').Value = '  -4.05%  '

# Row 12
$ws.Range('D12').Value = '2.290.11'
$ws.Range('E12').Value = '  -1.56%  '

# Row 13
$ws.Range('D13').Value = '''13.90'
$ws.Range('E13').Value = '  -7.87%  '

# Row 14
$ws.Range('D14').Value = '''20.98'
$ws.Range('E14').Value = '  -3.16%  '

# Row 15
$ws.Range('D15').Value = '''0.758'
$ws.Range('E15').Value = '  -9.25%  '

# Row 16
$ws.Range('D16').Value = '''5.07'
$ws.Range('E16').Value = '  -6.39%  '

# Row 17
$ws.Range('D17').Value = '1.976.40'
$ws.Range('E17').Value = '  -3.49%  '

# Row 18
$ws.Range('D18').Value = '36.971.05'
$ws.Range('E18').Value = '  -0.99%  '

# Row 19
$ws.Range('D19').Value = '''68.71'
$ws.Range('E19').Value = '  -2.23%  '

# Row 20
$ws.Range('D20').Value = '0.0₃0811'
$ws.Range('E20').Value = '  -5.92%  '

# Row 21
$ws.Range('D21').Value = '''229.90'
$ws.Range('E21').Value = '  +0.13%  '

# Row 22
$ws.Range('D22').Value = '''4.97'
$ws.Range('E22').Value = '  -6.00%  '

# Row 23
$ws.Range('D23').Value = '''0.997'
$ws.Range('E23').Value = '  -0.35%  '

# Row 24
$ws.Range('D24').Value = '''2.48'
$ws.Range('E24').Value = '  -6.44%  '

# Row 25
$ws.Range('D25').Value = '''2.35'
$ws.Range('E25').Value = '  -0.46%  '

# Row 26
$ws.Range('D26').Value = '''163.31'
$ws.Range('E26').Value = '  -0.31%  '

# Row 27
$ws.Range('D27').Value = '''8.65'
$ws.Range('E27').Value = '  -6.86%  '

# Row 28
$ws.Range('D28').Value = '''0.128'
$ws.Range('E28').Value = '  -7.37%  '

# Row 29
$ws.Range('D29').Value = '''19.07'
$ws.Range('E29').Value = '  -4.81%  '

# Row 30
$ws.Range('D30').Value = '''1.27'
$ws.Range('E30').Value = '  -7.10%  '

# Row 31
$ws.Range('D31').Value = '''0.117'
$ws.Range('E31').Value = '  -3.55%  '

# Row 32
$ws.Range('D32').Value = '''4.45'
$ws.Range('E32').Value = '  -7.43%  '

# Row 33
$ws.Range('D33').Value = '''0.0612'
$ws.Range('E33').Value = '  -9.27%  '

# Row 34
$ws.Range('D34').Value = '''4.27'
$ws.Range('E34').Value = '  -6.30%  '

# Row 35
$ws.Range('D35').Value = '''2.34'
$ws.Range('E35').Value = '  -6.34%  '

# Row 36
$ws.Range('E36').Value = '  -0.17%  '

# Row 37
$ws.Range('E37').Value = '  +0.41%  '

# Row 38
$ws.Range('D38').Value = '''3.35'
$ws.Range('E38').Value = '  -6.16%  '

# Row 39
$ws.Range('D39').Value = '''5.21'
$ws.Range('E39').Value = '  -3.68%  '

# Row 40
$ws.Range('D40').Value = '''3.03'
$ws.Range('E40').Value = '  +0.48%  '

# Row 41
$ws.Range('D41').Value = '1.423.37'
$ws.Range('E41').Value = '  +1.24%  '

# Row 42
$ws.Range('D42').Value = '''1.14'
$ws.Range('E42').Value = '  -3.88%  '

# Row 43
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '''0.0205'
$ws.Range('E43').Value = '  -6.17%  '

# Row 44
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').Value = '''0.0896'
$ws.Range('E44').Value = '  -8.48%  '

# Row 45
$ws.Range('D45').Value = '''15.44'
$ws.Range('E45').Value = '  -7.33%  '

# Row 46
$ws.Range('D46').Value = '''87.49'
$ws.Range('E46').Value = '  -4.73%  '

# Row 47
$ws.Range('D47').Value = '''1.00'
$ws.Range('E47').Value = '  -4.93%  '

# Row 48
$ws.Range('D48').Value = '''2.88'
$ws.Range('E48').Value = '  -0.22%  '

# Row 49
$ws.Range('D49').Value = '''6.68'
$ws.Range('E49').Value = '  -10.98%  '

# Row 50
$ws.Range('D50').Value = '2.184.72'
$ws.Range('E50').Value = '  -1.41%  '

# Row 51
$ws.Range('D51').Value = '''1.86'
$ws.Range('E51').Value = '  -10.83%  '
